$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column F (gene_type / Gene) and G (sourcespecies / Acetobacter xylinum)
# Order matters for shared string table indices: Acetobacter xylinum must become
# index 9 and "sourcespecies " index 10, so populate G2 before G1.
$ws.Range("F1").Value = "gene_type"
$ws.Range("F2").Value = "Gene"
$ws.Range("G2").Value = "Acetobacter xylinum"
$ws.Range("G1").Value = "sourcespecies "

# Update view: scroll so column B is leftmost, select C6
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C6").Select()

$excel.Calculate()
